# separando o comando de salvar dados da interface
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Programacao": append rows 40-42 (copies of the pattern
# used by the existing rows), dimension grows from K39 to K42.
# ---------------------------------------------------------------
$wsProg = $wb.Worksheets.Item("Programacao")

for ($r = 40; $r -le 42; $r++) {
    $wsProg.Cells.Item($r, 1).Value = "das"
    $wsProg.Cells.Item($r, 2).Value = "das"
    $wsProg.Cells.Item($r, 3).Value = "das"
    $wsProg.Cells.Item($r, 4).Value = "das"
    $wsProg.Cells.Item($r, 5).Value = "das"
    $wsProg.Cells.Item($r, 6).Value = "NORSAL"
    $wsProg.Cells.Item($r, 7).Value = 516
    $wsProg.Cells.Item($r, 8).Value = "das"
    $wsProg.Cells.Item($r, 9).Value = "das"
    $wsProg.Cells.Item($r, 10).Value = "SAL REFINADO Selecione uma opção"
    $wsProg.Cells.Item($r, 11).Value = "das"
}

# ---------------------------------------------------------------
# Sheet "Planilha": append rows 62-67, dimension grows from N61 to N67.
# ---------------------------------------------------------------
$wsPlan = $wb.Worksheets.Item("Planilha")

$planRows = @(
    @{ A="ENTRADA"; B="das"; C="das"; D="das"; E="das"; F="SAL REFINADO"; G="Selecione uma opção"; H="NORSAL"; I="das"; J="das"; K=150561; L="ddas"; M="das"; N=516 },
    @{ A="ENTRADA"; B="das"; C="das"; D="das"; E="das"; F="SAL REFINADO"; G="Selecione uma opção"; H="NORSAL"; I="das"; J="das"; K=150561; L="ddas"; M="das"; N=516 },
    @{ A="ENTRADA"; B="das"; C="das"; D="das"; E="das"; F="SAL REFINADO"; G="Selecione uma opção"; H="NORSAL"; I="das"; J="das"; K=5156;   L="das";  M="das"; N=356 },
    @{ A="ENTRADA"; B="das"; C="das"; D="das"; E="das"; F="SAL REFINADO"; G="Selecione uma opção"; H="NORSAL"; I="das"; J="das"; K=150561; L="ddas"; M="das"; N=516 },
    @{ A="ENTRADA"; B="das"; C="das"; D="das"; E="das"; F="SAL REFINADO"; G="Selecione uma opção"; H="NORSAL"; I="das"; J="das"; K=5156;   L="das";  M="das"; N=356 },
    @{ A="ENTRADA"; B="das"; C="das"; D="das"; E="das"; F="SAL REFINADO"; G="Selecione uma opção"; H="NORSAL"; I="das"; J="das"; K=86415;  L="das";  M="das"; N=6541 }
)

$r = 62
foreach ($row in $planRows) {
    $wsPlan.Cells.Item($r, 1).Value = $row.A
    $wsPlan.Cells.Item($r, 2).Value = $row.B
    $wsPlan.Cells.Item($r, 3).Value = $row.C
    $wsPlan.Cells.Item($r, 4).Value = $row.D
    $wsPlan.Cells.Item($r, 5).Value = $row.E
    $wsPlan.Cells.Item($r, 6).Value = $row.F
    $wsPlan.Cells.Item($r, 7).Value = $row.G
    $wsPlan.Cells.Item($r, 8).Value = $row.H
    $wsPlan.Cells.Item($r, 9).Value = $row.I
    $wsPlan.Cells.Item($r, 10).Value = $row.J
    $wsPlan.Cells.Item($r, 11).Value = $row.K
    $wsPlan.Cells.Item($r, 12).Value = $row.L
    $wsPlan.Cells.Item($r, 13).Value = $row.M
    $wsPlan.Cells.Item($r, 14).Value = $row.N
    $r++
}

# ---------------------------------------------------------------
# Sheet "Descarga do Sal": update several existing cells.
# ---------------------------------------------------------------
$wsDesc = $wb.Worksheets.Item("Descarga do Sal")

$wsDesc.Range("D8").Value = "das"
$wsDesc.Range("M18").Value = "Selecione uma opção"
$wsDesc.Range("P20").Value = 7413
$wsDesc.Range("D28").Value = "ddas"
$wsDesc.Range("K28").Value = 516
$wsDesc.Range("O28").Value = 150561
$wsDesc.Range("D30").Value = "das"
$wsDesc.Range("H30").Value = "das"
$wsDesc.Range("K30").Value = 356
$wsDesc.Range("O30").Value = 5156
$wsDesc.Range("D32").Value = "das"
$wsDesc.Range("H32").Value = "das"
$wsDesc.Range("K32").Value = 6541
$wsDesc.Range("O32").Value = 86415
